$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-9
# from 2023-09-14 (45183) to 2023-09-15 (45184)
$ws.Range("C2:C9").Value = 45184
